$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the MARCUS / 002606448 / 450000 row (was row 2, right under the header).
$ws.Rows(2).Delete()

# 2) ANILSON's Saldo (004385806) changes from 179359.28 to 170129.36.
#    After the deletion above, ANILSON is now on row 2.
$ws.Range("C2").Value = 170129.36

# 3) Insert a new row for LEVI (005206566 / 50000) right before BRASFORT (004352384),
#    which is now on row 4.
$ws.Rows(4).Insert()
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "005206566"
$ws.Range("B4").Value = "LEVI"
$ws.Range("C4").Value = 50000

# 4) Remove the JOSE (004432935 / 10000) and CESAR (004207278 / 5780) rows, which
#    now sit right after ANA (004479287 / 20000) on rows 9 and 10.
$ws.Rows(10).Delete()
$ws.Rows(9).Delete()
